# Update crypto price/volume data as scraped on Mon Mar 18 02:25:21 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.634.96"
$ws.Range("E2").Value = "  +2.17%  "

# Row 3
$ws.Range("D3").Value = "3.599.86"
$ws.Range("E3").Value = "  +1.01%  "

# Row 4
$ws.Range("E4").Value = "  +0.20%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "199.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.70%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "559.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.45%  "

# Row 7
$ws.Range("D7").Value = "3.593.31"
$ws.Range("E7").Value = "  +1.09%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.614"
$ws.Range("D8").Style = "Normal"

# Row 9
$ws.Range("E9").Value = "  -0.27%  "

# Row 10
$ws.Range("E10").Value = "  +0.41%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "58.92"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +10.36%  "

# Row 12
$ws.Range("E12").Value = "  +4.35%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000287"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +12.84%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.98"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.71%  "

# Row 15
$ws.Range("D15").Value = "4.192.62"
$ws.Range("E15").Value = "  +1.55%  "

# Row 16
$ws.Range("D16").Value = "3.605.56"
$ws.Range("E16").Value = "  +1.16%  "

# Row 17
$ws.Range("E17").Value = "  +0.50%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.88"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.49%  "

# Row 19
$ws.Range("D19").Value = "67.601.32"
$ws.Range("E19").Value = "  +2.44%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.62%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.23%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "398.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.92%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +19.73%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.13"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.73%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.02"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.30%  "

# Row 26
$ws.Range("E26").Value = "  +2.83%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.33%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.26%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +8.97%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.51"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +22.83%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.42%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.44"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.23%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "664.77"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.67%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "12.18"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.60%  "

# Row 35
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "63.61"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.45%  "

# Row 36
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.113"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.18%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "42.32"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.79%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.431"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +13.27%  "

# Row 39
$ws.Range("E39").Value = "  -0.03%  "

# Row 40
$ws.Range("D40").Value = "0.0₃0771"
$ws.Range("E40").Value = "  +3.65%  "

# Row 41
$ws.Range("E41").Value = "  +14.76%  "

# Row 42
$ws.Range("D42").Value = "3.248.65"
$ws.Range("E42").Value = "  +10.40%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.136"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.17%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.81"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +14.88%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.17%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +30.52%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0416"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.98%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.72"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +10.73%  "

# Row 49
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.12"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.84%  "

# Row 50
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.130"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.21%  "

# Row 51
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.97%  "
